$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Update changed cell values
$ws.Range("J1").Value = 33.03995990753174
$ws.Range("J2").Value = 43.27023649215698
$ws.Range("B3").Value = 2084
$ws.Range("J3").Value = 39.24678134918213
$ws.Range("B4").Value = 2586
$ws.Range("D4").Value = 2563
$ws.Range("E4").Value = 22
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = 99.68883702839362
$ws.Range("H4").Value = 99.14893617021276
$ws.Range("I4").Value = 0.01166407465007776
$ws.Range("J4").Value = 35.72079944610596
$ws.Range("B5").Value = 2021
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 15
$ws.Range("G5").Value = 99.25962487660415
$ws.Range("H5").Value = 99.55445544554455
$ws.Range("I5").Value = 0.01184015786877158
$ws.Range("J5").Value = 43.63012886047363
$ws.Range("B6").Value = 1760
$ws.Range("D6").Value = 1756
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 99.65947786606129
$ws.Range("H6").Value = 99.82944855031268
$ws.Range("I6").Value = 0.005104934770277935
$ws.Range("J6").Value = 35.04046368598938
$ws.Range("B7").Value = 2527
$ws.Range("D7").Value = 2526
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 99.80244962465429
$ws.Range("I7").Value = 0.001974723538704581
$ws.Range("J7").Value = 39.14115214347839
$ws.Range("J8").Value = 33.55628561973572
$ws.Range("J9").Value = 37.80549788475037
$ws.Range("B10").Value = 1809
$ws.Range("D10").Value = 1793
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 99.94425863991081
$ws.Range("H10").Value = 99.17035398230088
$ws.Range("I10").Value = 0.008913649025069638
$ws.Range("J10").Value = 37.77962040901184
$ws.Range("B11").Value = 1870
$ws.Range("D11").Value = 1868
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 99.46751863684771
$ws.Range("H11").Value = 99.94649545211342
$ws.Range("I11").Value = 0.005854177754124534
$ws.Range("J11").Value = 28.80000448226929
$ws.Range("J12").Value = 32.71163868904114
$ws.Range("J13").Value = 32.22668743133545
$ws.Range("J14").Value = 31.75705552101135
$ws.Range("B15").Value = 2280
$ws.Range("E15").Value = 2
$ws.Range("H15").Value = 99.91224221149628
$ws.Range("I15").Value = 0.000877963125548727
$ws.Range("J15").Value = 33.2398419380188
$ws.Range("B16").Value = 2000
$ws.Range("E16").Value = 13
$ws.Range("H16").Value = 99.3496748374187
$ws.Range("I16").Value = 0.006542526421741319
$ws.Range("J16").Value = 37.01126217842102
$ws.Range("J17").Value = 36.77881526947021
$ws.Range("J18").Value = 35.81173610687256
$ws.Range("B19").Value = 1518
$ws.Range("E19").Value = 0
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 35.65586733818054
$ws.Range("J20").Value = 28.28019499778748
$ws.Range("B21").Value = 2600
$ws.Range("E21").Value = 2
$ws.Range("H21").Value = 99.92304732589457
$ws.Range("I21").Value = 0.001922337562475971
$ws.Range("J21").Value = 35.35386109352112
$ws.Range("J22").Value = 41.25642824172974
$ws.Range("B23").Value = 2048
$ws.Range("D23").Value = 2045
$ws.Range("F23").Value = 90
$ws.Range("G23").Value = 95.78454332552693
$ws.Range("H23").Value = 99.90229604298975
$ws.Range("I23").Value = 0.04307116104868914
$ws.Range("J23").Value = 32.55426836013794
$ws.Range("B24").Value = 2949
$ws.Range("D24").Value = 2934
$ws.Range("E24").Value = 14
$ws.Range("F24").Value = 45
$ws.Range("G24").Value = 98.48942598187311
$ws.Range("H24").Value = 99.52510176390773
$ws.Range("I24").Value = 0.01979865771812081
$ws.Range("J24").Value = 41.49635148048401
$ws.Range("B25").Value = 2645
$ws.Range("D25").Value = 2642
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 99.51035781544256
$ws.Range("H25").Value = 99.92435703479576
$ws.Range("I25").Value = 0.005647590361445783
$ws.Range("J25").Value = 35.89228439331055
$ws.Range("B26").Value = 1852
$ws.Range("D26").Value = 1848
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 99.40828402366864
$ws.Range("H26").Value = 99.83792544570503
$ws.Range("I26").Value = 0.007526881720430108
$ws.Range("J26").Value = 37.92574667930603
$ws.Range("D27").Value = 2941
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 99.55991875423155
$ws.Range("H27").Value = 99.89809782608695
$ws.Range("I27").Value = 0.005414551607445008
$ws.Range("J27").Value = 40.86296820640564
$ws.Range("B28").Value = 3006
$ws.Range("E28").Value = 1
$ws.Range("H28").Value = 99.96672212978369
$ws.Range("I28").Value = 0.0003327787021630616
$ws.Range("J28").Value = 36.12000465393066
$ws.Range("B29").Value = 2613
$ws.Range("D29").Value = 2609
$ws.Range("E29").Value = 3
$ws.Range("F29").Value = 40
$ws.Range("G29").Value = 98.48999622499056
$ws.Range("H29").Value = 99.88514548238898
$ws.Range("I29").Value = 0.01622641509433962
$ws.Range("J29").Value = 36.86920571327209
$ws.Range("J30").Value = 35.20494389533997
$ws.Range("B31").Value = 3247
$ws.Range("D31").Value = 3246
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 99.87692307692308
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 0.001230390649031067
$ws.Range("J31").Value = 36.66146159172058
$ws.Range("B32").Value = 2259
$ws.Range("E32").Value = 1
$ws.Range("H32").Value = 99.95571302037202
$ws.Range("I32").Value = 0.002210433244916004
$ws.Range("J32").Value = 38.30744004249573
$ws.Range("J33").Value = 37.26217818260193
$ws.Range("J34").Value = 37.99660730361938
$ws.Range("J35").Value = 44.58633685112
$ws.Range("D36").Value = 2417
$ws.Range("E36").Value = 6
$ws.Range("F36").Value = 9
$ws.Range("G36").Value = 99.6290189612531
$ws.Range("H36").Value = 99.75237309120925
$ws.Range("I36").Value = 0.006180469715698393
$ws.Range("J36").Value = 36.40898275375366
$ws.Range("B37").Value = 2451
$ws.Range("D37").Value = 2449
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 33
$ws.Range("G37").Value = 98.67042707493957
$ws.Range("H37").Value = 99.95918367346938
$ws.Range("I37").Value = 0.01369311316955296
$ws.Range("J37").Value = 40.3929455280304
$ws.Range("B38").Value = 2589
$ws.Range("D38").Value = 2588
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 16
$ws.Range("G38").Value = 99.38556067588326
$ws.Range("H38").Value = 100
$ws.Range("I38").Value = 0.006142034548944338
$ws.Range("J38").Value = 31.97047209739685
$ws.Range("D39").Value = 2047
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 5
$ws.Range("G39").Value = 99.75633528265107
$ws.Range("H39").Value = 99.65920155793573
$ws.Range("I39").Value = 0.005845104724792985
$ws.Range("J39").Value = 37.92467474937439
$ws.Range("J40").Value = 39.04234647750854
$ws.Range("J41").Value = 35.44862580299377
$ws.Range("B42").Value = 1782
$ws.Range("D42").Value = 1777
$ws.Range("E42").Value = 4
$ws.Range("F42").Value = 2
$ws.Range("G42").Value = 99.8875772906127
$ws.Range("H42").Value = 99.77540707467715
$ws.Range("I42").Value = 0.003370786516853933
$ws.Range("J42").Value = 34.59914135932922
$ws.Range("D43").Value = 3076
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 2
$ws.Range("G43").Value = 99.93502274204029
$ws.Range("H43").Value = 100
$ws.Range("I43").Value = 0.0006495615459564793
$ws.Range("J43").Value = 43.10721158981323
$ws.Range("J44").Value = 36.31776738166809

# Rename sheet 2
$ws.Name = "shiccsd"
